# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" quarter:
#   - total/summary sheet ("总计") gets a new row 2 with the 2022-Q4 rollup,
#     existing quarter rows shift down one row
#   - a brand new detail worksheet named "2022-Q4" is added (copied from the
#     "2022-Q3" sheet so it keeps the same look/formatting), placed right
#     after "总计" and before "2022-Q3"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet ("总计"): shift existing quarter rows down one row, then
#    write the new 2022-Q4 rollup into row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Walk from the bottom up so we never overwrite data before reading it.
$total.Range("B7").Value = $total.Range("B6").Text
$total.Range("C7").Value = $total.Range("C6").Value
$total.Range("D7").Value = $total.Range("D6").Value

$total.Range("B6").Value = $total.Range("B5").Text
$total.Range("C6").Value = $total.Range("C5").Value
$total.Range("D6").Value = $total.Range("D5").Value

$total.Range("B5").Value = $total.Range("B4").Text
$total.Range("C5").Value = $total.Range("C4").Value
$total.Range("D5").Value = $total.Range("D4").Value

$total.Range("B4").Value = $total.Range("B3").Text
$total.Range("C4").Value = $total.Range("C3").Value
$total.Range("D4").Value = $total.Range("D3").Value

$total.Range("B3").Value = $total.Range("B2").Text
$total.Range("C3").Value = $total.Range("C2").Value
$total.Range("D3").Value = $total.Range("D2").Value

# New 2022-Q4 rollup row.
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.27

# Re-number the helper index column (A) 0..5 and make sure every row in the
# now 6-row table carries the same style as the rest of column A / header.
$total.Range("B1").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# 2. New detail worksheet "2022-Q4": copy the "2022-Q3" sheet (so headers,
#    column styles etc. match) and place it right before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Overwrite the (copied) data rows with the 2022-Q4 fund holdings.
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'013442"
$q4.Range("C2").Value = "建信中证1000指数增强E"
$q4.Range("D2").Value = "'9.52"
$q4.Range("E2").Value = "'86.80"
$q4.Range("F2").Value = "'1.11"
$q4.Range("G2").Value = "'0.1057"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'006165"
$q4.Range("C3").Value = "建信中证1000指数增强A"
$q4.Range("D3").Value = "'7.20"
$q4.Range("E3").Value = "'86.80"
$q4.Range("F3").Value = "'1.11"
$q4.Range("G3").Value = "'0.0799"
$q4.Range("H3").Value = 9

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'000965"
$q4.Range("C4").Value = "汇丰晋信新动力混合"
$q4.Range("D4").Value = "'0.95"
$q4.Range("E4").Value = "'91.04"
$q4.Range("F4").Value = "'4.14"
$q4.Range("G4").Value = "'0.0393"
$q4.Range("H4").Value = 4

# Rows 5 & 6 are new - give column A the same style the rest of column A uses
# before filling in the values.
$q4.Range("A2").Copy()
$q4.Range("A5:A6").PasteSpecial(-4122)

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'006166"
$q4.Range("C5").Value = "建信中证1000指数增强C"
$q4.Range("D5").Value = "'2.21"
$q4.Range("E5").Value = "'86.80"
$q4.Range("F5").Value = "'1.11"
$q4.Range("G5").Value = "'0.0245"
$q4.Range("H5").Value = 9

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'540004"
$q4.Range("C6").Value = "汇丰晋信2026周期混合"
$q4.Range("D6").Value = "'1.08"
$q4.Range("E6").Value = "'23.97"
$q4.Range("F6").Value = "'1.99"
$q4.Range("G6").Value = "'0.0215"
$q4.Range("H6").Value = 2
